$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "Globo"
$ws.Range("B13").Value = "RJ TV 1"
$ws.Range("C13").Value = "Trânsito"
$ws.Range("D13").Value = "2025-04-01T12:52"
$ws.Range("E13").Value = "Neutro"
$ws.Range("F13").Value = "Caminhão invade o calçadão de Campos. Motorista teria errado o caminho e acabou subindo no calçadão do Centro. Repórter *ao vivo* do local. Imagens no Boulevard Francisco de Paula Carneiro, no Centro. Sem energia. Equipe no local"
